$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" column (E) so the list of periods is now shown in
# reverse-chronological order (most recent period first).
$ws.Range("E16").Value = "2104"
$ws.Range("E17").Value = "2103"
$ws.Range("E18").Value = "2102"
$ws.Range("E19").Value = "2101"
$ws.Range("E20").Value = "2012"
$ws.Range("E21").Value = "2011"
$ws.Range("E22").Value = "2010"
$ws.Range("E23").Value = "2009"
$ws.Range("E24").Value = "2008"
$ws.Range("E25").Value = "2007"
$ws.Range("E26").Value = "2006"
$ws.Range("E27").Value = "2005"
$ws.Range("E28").Value = "2004"
$ws.Range("E29").Value = "2003"
$ws.Range("E30").Value = "2002"
$ws.Range("E31").Value = "2001"
$ws.Range("E32").Value = "1912"

# Swap the "Valor Mora" amounts between the first and last data rows so the
# figures stay aligned with their corresponding periods.
$ws.Range("F16").Value = 29398
$ws.Range("F32").Value = 33920
